$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected Price/Volume columns to text so values like "307.78" or "-6.21%"
# are stored as literal strings (matching original inlineStr cells) rather than
# being auto-converted to numbers/percentages by Excel.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '307.78'
$ws.Range("E2").Value = '-6.21%'
$ws.Range("D3").Value = '39.89'
$ws.Range("E3").Value = '-10.12%'
$ws.Range("D4").Value = '5.068'
$ws.Range("E4").Value = '-5.24%'
$ws.Range("D5").Value = '0.07785'
$ws.Range("E5").Value = '-6.94%'
$ws.Range("D6").Value = '4.323'
$ws.Range("E6").Value = '-2.10%'
$ws.Range("D7").Value = '1.655'
$ws.Range("E7").Value = '-14.66%'
$ws.Range("D8").Value = '0.9175'
$ws.Range("E8").Value = '-5.63%'
$ws.Range("D9").Value = '0.09933'
$ws.Range("E9").Value = '-11.95%'
$ws.Range("D10").Value = '0.1738'
$ws.Range("E10").Value = '-8.45%'
$ws.Range("D11").Value = '0.08934'
$ws.Range("E11").Value = '-7.74%'
$ws.Range("D12").Value = '0.04393'
$ws.Range("E12").Value = '-4.58%'
$ws.Range("D13").Value = '7.053'
$ws.Range("E13").Value = '-15.04%'
$ws.Range("D14").Value = '0.1060'
$ws.Range("E14").Value = '-0.07%'
$ws.Range("D15").Value = '0.001255'
$ws.Range("E15").Value = '-3.87%'
$ws.Range("D16").Value = '0.005654'
$ws.Range("E16").Value = '-2.43%'
$ws.Range("D17").Value = '3.366'
$ws.Range("E17").Value = '0.16%'
$ws.Range("E19").Value = '0.28%'
$ws.Range("D20").Value = '0.1365'
$ws.Range("E20").Value = '-1.65%'
$ws.Range("E21").Value = '0.22%'
$ws.Range("D22").Value = '0.04145'
$ws.Range("E22").Value = '-0.84%'
$ws.Range("D23").Value = '0.001206'
$ws.Range("E23").Value = '-2.85%'
$ws.Range("E24").Value = '-7.71%'
$ws.Range("D25").Value = '0.0001224'
$ws.Range("E25").Value = '-5.82%'
$ws.Range("D26").Value = '0.0002994'
$ws.Range("E26").Value = '0.49%'
$ws.Range("D38").Value = '0.02389'
$ws.Range("E38").Value = '-12.07%'
$ws.Range("E39").Value = '-7.67%'
$ws.Range("D40").Value = '0.007989'
$ws.Range("E40").Value = '2.13%'
$ws.Range("D41").Value = '0.1325'
$ws.Range("E41").Value = '-6.31%'
$ws.Range("D42").Value = '0.007106'
$ws.Range("E42").Value = '-3.16%'
$ws.Range("D43").Value = '0.002018'
$ws.Range("E43").Value = '-1.06%'
$ws.Range("D44").Value = '0.008053'
$ws.Range("E44").Value = '-7.50%'
$ws.Range("D45").Value = '0.3337'
$ws.Range("E45").Value = '-5.04%'
$ws.Range("D46").Value = '0.00006726'
$ws.Range("E46").Value = '-2.66%'
$ws.Range("E47").Value = '0.35%'
$ws.Range("D48").Value = '0.003422'
$ws.Range("E48").Value = '-1.82%'
$ws.Range("D49").Value = '0.004120'
$ws.Range("E49").Value = '16.72%'
$ws.Range("D50").Value = '0.00002107'
$ws.Range("E50").Value = '0.35%'
$ws.Range("D51").Value = '0.0002007'
$ws.Range("E51").Value = '0.35%'

# Restore the default (unstyled) cell style now that the values are committed as text,
# so no stray number-format style is left attached to these cells.
$priceVolRange.Style = "Normal"
